# Update gh-pages to output generated at 456a3b4
#
# Updates the "想去人数" (F column) interest counters - and one cover-image
# URL - across the 展览 (sheet1), 本地生活 (sheet3) and 全部类型 (sheet4)
# worksheets to reflect freshly scraped counts.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")

$ws.Range("F2").Value2 = 670
$ws.Range("F3").Value2 = 47
$ws.Range("F4").Value2 = 1986
$ws.Range("F5").Value2 = 5765
$ws.Range("F6").Value2 = 1612
$ws.Range("F8").Value2 = 3259
$ws.Range("F11").Value2 = 1370
$ws.Range("F12").Value2 = 4550
$ws.Range("F13").Value2 = 1086
$ws.Range("F14").Value2 = 1716
$ws.Range("F15").Value2 = 2606
$ws.Range("F16").Value2 = 5
$ws.Range("F17").Value2 = 49
$ws.Range("F21").Value2 = 1026
$ws.Range("F23").Value2 = 83
$ws.Range("F24").Value2 = 16
$ws.Range("I24").Value = "//i1.hdslb.com/bfs/openplatform/202404/YYAGMoXP1714288325893.jpeg"
$ws.Range("F29").Value2 = 1122
$ws.Range("F30").Value2 = 411
$ws.Range("F31").Value2 = 88
$ws.Range("F32").Value2 = 205
$ws.Range("F33").Value2 = 385
$ws.Range("F34").Value2 = 920
$ws.Range("F36").Value2 = 1746
$ws.Range("F37").Value2 = 2255
$ws.Range("F38").Value2 = 1051
$ws.Range("F40").Value2 = 273
$ws.Range("F42").Value2 = 377
$ws.Range("F43").Value2 = 39
$ws.Range("F44").Value2 = 672
$ws.Range("F45").Value2 = 31
$ws.Range("F46").Value2 = 444
$ws.Range("F47").Value2 = 407
$ws.Range("F49").Value2 = 147

# ---------------------------------------------------------------------
# Sheet "本地生活"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Range("F2").Value2 = 782

# ---------------------------------------------------------------------
# Sheet "全部类型"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value2 = 782
$ws4.Range("F3").Value2 = 670
$ws4.Range("F4").Value2 = 47
$ws4.Range("F5").Value2 = 1986
$ws4.Range("F6").Value2 = 5765
$ws4.Range("F7").Value2 = 1612
$ws4.Range("F9").Value2 = 3259
$ws4.Range("F11").Value2 = 1370
$ws4.Range("F12").Value2 = 4550
$ws4.Range("F13").Value2 = 1086
$ws4.Range("F14").Value2 = 1716
$ws4.Range("F15").Value2 = 5
$ws4.Range("F17").Value2 = 49
$ws4.Range("F23").Value2 = 1026
$ws4.Range("F25").Value2 = 83
$ws4.Range("F29").Value2 = 1122
$ws4.Range("F30").Value2 = 411
$ws4.Range("F31").Value2 = 88
$ws4.Range("F32").Value2 = 205
$ws4.Range("F33").Value2 = 920
$ws4.Range("F34").Value2 = 1746
$ws4.Range("F35").Value2 = 2255
$ws4.Range("F36").Value2 = 1051
$ws4.Range("F40").Value2 = 273
$ws4.Range("F42").Value2 = 377
$ws4.Range("F43").Value2 = 672
$ws4.Range("F44").Value2 = 444
$ws4.Range("F45").Value2 = 407
$ws4.Range("F48").Value2 = 147
